$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(45922, 5049.38273576318, 4741.91257002738, 6432, 6527.061986, -8.85034082232494),
  @(45923, 5046.485857625, 4955.71415683252, 2952, 6523.437389, 145.027737008646),
  @(45924, 5032.92426258885, 5178.41364632139, 2952, 6505.013449, 154.10428469719),
  @(45925, 5053.65330497704, 5517.99142251688, 2952, 6531.364571, 168.487612022493),
  @(45926, 5223.86588860041, 4930.86458394508, 2952, 6756.427329, 146.309417681028),
  @(45927, 1198.04912219854, 2895.30236148839, 2952, 2322.416874, 44.4862547204102),
  @(45928, 1070.29084929054, 2779.28377354922, 2952, 2183.81959, 39.2005214274454),
  @(45929, 5457.44872542722, 5419.69389657574, 2952, 7030.796743, 168.376746422855),
  @(45930, 5457.44872542722, 5270.86550023517, 2952, 7030.796743, 162.175563241998),
  @(45931, 4260.05383201679, 4487.26418315459, 3692, 6097.879545, 109.712079005742),
  @(45932, 4260.05383201679, 4727.14983180789, 3692, 6097.879545, 119.707314366296),
  @(45933, 4260.05383201679, 4119.50232522947, 3692, 6097.879545, 94.3886682588618),
  @(45934, 805.573379841386, 1966.81059557387, 3692, 2301.949365, -9.53389246947986),
  @(45935, 709.592596030573, 1826.70918032675, 3692, 2197.537945, -15.7227279459928)
)

for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Length; $j++) {
    $ws.Cells.Item($i + 2, $j + 1).Value = $row[$j]
  }
}

Write-Host "done"